$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.05093278012031993
$ws.Range("C2").Value = 0.9985021538304686
$ws.Range("D2").Value = 0.16661417163493
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 AdaBoostRegressor(learning_rate=0.5, n_estimators=100))])"
$ws.Range("G2").Value = 0.1248244242667473
$ws.Range("H2").Value = 0.99
